$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32
$ws.Range("C32").Value = "Verify that placeholder and tooltip are present or not."
$ws.Range("B32").Value = "TC-012"

# Row 33
$ws.Range("B33").Value = "TC-013"
$ws.Range("C33").Value = "Verify that the login screen is responsive according to different screen sizes. (Responsive Testing)."

# Row 34 / Row 35 (test case numbers entered first, descriptions after)
$ws.Range("B34").Value = "TC-014"
$ws.Range("B35").Value = "TC-015"
$ws.Range("C34").Value = "Verify that user is able to input username and password."
$ws.Range("C35").Value = "Verify that password is masked/encrypted."

# Row 36
$ws.Range("C36").Value = "Verify that password should be visible when eye button is clicked and vice-versa."
$ws.Range("B36").Value = "TC-016"

# Rows 37-39 (Test case numbers only, no description yet)
$ws.Range("B37").Value = "TC-017"
$ws.Range("B38").Value = "TC-018"
$ws.Range("B39").Value = "TC-019"

# Row heights to match wrapped multi-line text (rows 33 & 36 wrap to 2 lines)
$ws.Rows("33").RowHeight = 30
$ws.Rows("36").RowHeight = 30

# Restore the sheet view scroll / selection state after data entry
$ws.Range("C37").Select() | Out-Null
